# Update recomputed radiomics feature values in row 2 (M013 patient record)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = [double]"0.60459186836377932"
$ws.Range("M2").Value = [double]"0.88623629890251221"
$ws.Range("O2").Value = [double]"0.00042768091514865038"
$ws.Range("P2").Value = [double]"218.29863997495968"
$ws.Range("Q2").Value = [double]"11.730759074140661"
$ws.Range("R2").Value = [double]"0.18266217861447268"
$ws.Range("S2").Value = [double]"0.6395127083510973"
$ws.Range("T2").Value = [double]"0.0006356857814803895"
$ws.Range("U2").Value = [double]"0.0046551591305755338"
$ws.Range("V2").Value = [double]"10.874411901789523"
$ws.Range("W2").Value = [double]"1902.173312359984"
$ws.Range("X2").Value = [double]"0.00023773443861522901"
$ws.Range("Y2").Value = [double]"0.59477883102754692"
$ws.Range("Z2").Value = [double]"1.5187467909051835"
$ws.Range("AA2").Value = [double]"18865.624796921511"
$ws.Range("AB2").Value = [double]"0.55757708307883747"
$ws.Range("AC2").Value = [double]"0.97402629503249349"
$ws.Range("AD2").Value = [double]"1.1132563672316449"
$ws.Range("AE2").Value = [double]"0.016597390152423088"
$ws.Range("AF2").Value = [double]"0.93362721847064445"
$ws.Range("AG2").Value = [double]"0.96503658621998134"
$ws.Range("AH2").Value = [double]"0.0064094772022125568"
$ws.Range("AI2").Value = [double]"1993.350283531724"
$ws.Range("AJ2").Value = [double]"0.0058967952618072372"
$ws.Range("AK2").Value = [double]"1941.4358767275767"
$ws.Range("AL2").Value = [double]"0.0094742000446225924"
$ws.Range("AM2").Value = [double]"2226.742704720014"
$ws.Range("AN2").Value = [double]"0.012560704338504342"
$ws.Range("AO2").Value = [double]"7.3608860875791276e-06"
$ws.Range("AP2").Value = [double]"0.76834424205309049"
$ws.Range("AQ2").Value = [double]"5.1400460777278392"
$ws.Range("AR2").Value = [double]"0.016044629285140184"
$ws.Range("AS2").Value = [double]"0.5473611030695883"
$ws.Range("AT2").Value = [double]"0.62692681078777257"
$ws.Range("AU2").Value = [double]"0.0034550727779965591"
$ws.Range("AV2").Value = [double]"2038.3389734486536"
$ws.Range("AW2").Value = [double]"0.0024279394443928594"
$ws.Range("AX2").Value = [double]"1584.719424526837"
$ws.Range("AY2").Value = [double]"0.1747994814262562"
$ws.Range("AZ2").Value = [double]"11074.678815575784"
$ws.Range("BA2").Value = [double]"0.00089774126836235661"
$ws.Range("BB2").Value = [double]"6.1190775172951716e-07"
